$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
Write-Host "HandoutMaster: $hm"
